{"js": "// Reword the sentence in the Zimbabwe paragraph of section 1\n// (\"research motivation\"): the phrase\n//   \"have plagued Zimbabwe for decades,\"\n// becomes\n//   \"are among the economic and social problems plaguing Zimbabwe,\"\n// (rest of the paragraph is unchanged).\n\nconst body = context.document.body;\n\nconst searchText = \"have plagued Zimbabwe for decades,\";\nconst replacementText =\n  \"are among the economic and social problems plaguing Zimbabwe,\";\n\nconst results = body.search(searchText, {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace only the first (and only) occurrence, preserving the\n  // surrounding text / formatting of the rest of the paragraph.\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Reword the sentence in the Zimbabwe paragraph of section 1\n# (\"research motivation\"): the phrase\n#   \"have plagued Zimbabwe for decades,\"\n# becomes\n#   \"are among the economic and social problems plaguing Zimbabwe,\"\n# (rest of the paragraph is unchanged).\n\n$d = $word.ActiveDocument\n\n$searchText = \"have plagued Zimbabwe for decades,\"\n$replacementText = \"are among the economic and social problems plaguing Zimbabwe,\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\n    $searchText,      # FindText\n    $false,           # MatchCase\n    $false,           # MatchWholeWord\n    $false,           # MatchWildcards\n    $false,           # MatchSoundsLike\n    $false,           # MatchAllWordForms\n    $true,            # Forward\n    1,                # Wrap (wdFindContinue)\n    $false,           # Format\n    $replacementText, # ReplaceWith\n    2                 # Replace (wdReplaceOne)\n) | Out-Null\n"}
